$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new "REFRESH" translation row is inserted above the existing row 33
# ("ANSWER50QUESTIONS"), pushing all the following rows down by one.
$ws.Rows.Item(33).Insert()

# Fill in the new row: Key, Polish, English
$ws.Cells.Item(33, 1).Value = "REFRESH"
$ws.Cells.Item(33, 2).Value = "Odśwież"
$ws.Cells.Item(33, 3).Value = "Refresh"

# Match the selection state recorded in the workbook after the edit
$ws.Range("C33").Select()
